$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 21950.766
$ws.Range("J17").Value = 21950.766
$ws.Range("L17").Value = 65852.298
$ws.Range("N17").Value = -66188.298
$ws.Range("H43").Value = 2836.842
$ws.Range("I43").Value = 2755.5557
$ws.Range("J43").Value = 2910
$ws.Range("K43").Value = 2755.5557
$ws.Range("L43").Value = 2910
$ws.Range("M43").Value = -2686.5557
$ws.Range("N43").Value = -3048
$ws.Range("H131").Value = 2313.6
$ws.Range("I131").Value = 1226.7273
$ws.Range("J131").Value = 5302.5
$ws.Range("K131").Value = 3680.1819
$ws.Range("L131").Value = 15907.5
$ws.Range("M131").Value = 1359.8181
$ws.Range("N131").Value = -25987.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 5728.25
$ws.Range("I43").Value = 842
$ws.Range("J43").Value = 6172.4546
$ws.Range("K43").Value = 842
$ws.Range("L43").Value = 6172.4546
$ws.Range("M43").Value = -529
$ws.Range("N43").Value = -6798.4546
$ws.Range("H122").Value = 1867.5333
$ws.Range("I122").Value = 1649.9
$ws.Range("K122").Value = 4949.700000000001
$ws.Range("M122").Value = -2499.700000000001
$ws.Range("H123").Value = 49429
$ws.Range("J123").Value = 49429
$ws.Range("L123").Value = 49429
$ws.Range("N123").Value = -59229
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H129").Value = 249999.5
$ws.Range("J129").Value = 249999.5
$ws.Range("L129").Value = 249999.5
$ws.Range("N129").Value = -259999.5
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = 0
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = 0

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 200031800
$ws.Range("H27").Value = 200031800
$ws.Range("H31").Value = 4118585
$ws.Range("I31").Value = 1656.2693
$ws.Range("J31").Value = 6064769.5
$ws.Range("K31").Value = 1656.2693
$ws.Range("L31").Value = 6064769.5
$ws.Range("M31").Value = -1361.2693
$ws.Range("N31").Value = -6065359.5
$ws.Range("H34").Value = 4118585
$ws.Range("I34").Value = 1656.2693
$ws.Range("J34").Value = 6064769.5
$ws.Range("K34").Value = 1656.2693
$ws.Range("L34").Value = 6064769.5
$ws.Range("M34").Value = -1454.2693
$ws.Range("N34").Value = -6065173.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 25303.166
$ws.Range("I18").Value = 27576.182
$ws.Range("K18").Value = 82728.546
$ws.Range("M18").Value = -82559.546
$ws.Range("H68").Value = 1404.1562
$ws.Range("I68").Value = 1079.091
$ws.Range("K68").Value = 3237.273
$ws.Range("M68").Value = -2426.273
$ws.Range("H71").Value = 1404.1562
$ws.Range("I71").Value = 1079.091
$ws.Range("K71").Value = 9711.819
$ws.Range("M71").Value = -5655.819
$ws.Range("H140").Value = 3485.2307
$ws.Range("I140").Value = 899.6667
$ws.Range("J140").Value = 4260.9
$ws.Range("K140").Value = 2699.0001
$ws.Range("L140").Value = 12782.7
$ws.Range("M140").Value = 2480.9999
$ws.Range("N140").Value = -23142.7

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3592.5
$ws.Range("I122").Value = 5185
$ws.Range("K122").Value = 15555
$ws.Range("M122").Value = -13105
$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226
$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920
$ws.Range("H128").Value = 40001
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 40001
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 40001
$ws.Range("M128").ClearContents()
$ws.Range("N128").Value = -49961
$ws.Range("H129").Value = 50001
$ws.Range("J129").Value = 50001
$ws.Range("L129").Value = 50001
$ws.Range("N129").Value = -60001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 76926590
$ws.Range("I7").Value = 90911610
$ws.Range("J7").Value = 9002.5
$ws.Range("K7").Value = 90911610
$ws.Range("L7").Value = 9002.5
$ws.Range("M7").Value = -90911498
$ws.Range("N7").Value = -9226.5
$ws.Range("H40").Value = 5347
$ws.Range("I40").Value = 5056
$ws.Range("K40").Value = 5056
$ws.Range("M40").Value = -4920
$ws.Range("H94").Value = 100000
$ws.Range("J94").Value = 100000
$ws.Range("L94").Value = 100000
$ws.Range("N94").Value = -101352
$ws.Range("H126").Value = 76926590
$ws.Range("I126").Value = 90911610
$ws.Range("J126").Value = 9002.5
$ws.Range("K126").Value = 272734830
$ws.Range("L126").Value = 27007.5
$ws.Range("M126").Value = -272732360
$ws.Range("N126").Value = -31947.5
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 1554.8605
$ws.Range("I136").Value = 1178.3823
$ws.Range("J136").Value = 2977.111
$ws.Range("K136").Value = 3535.1469
$ws.Range("L136").Value = 8931.332999999999
$ws.Range("M136").Value = -985.1468999999997
$ws.Range("N136").Value = -14031.333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 799
$ws.Range("I122").Value = 799
$ws.Range("K122").Value = 2397
$ws.Range("M122").Value = 53
$ws.Range("H123").Value = 35000
$ws.Range("J123").Value = 35000
$ws.Range("L123").Value = 35000
$ws.Range("N123").Value = -44800
$ws.Range("H126").Value = 4903709.5
$ws.Range("I126").Value = 7354316
$ws.Range("K126").Value = 22062948
$ws.Range("M126").Value = -22060478
$ws.Range("H132").Value = 1967.6774
$ws.Range("I132").Value = 1223.9048
$ws.Range("J132").Value = 3529.6
$ws.Range("K132").Value = 3671.7144
$ws.Range("L132").Value = 10588.8
$ws.Range("M132").Value = -1141.7144
$ws.Range("N132").Value = -15648.8
